# Finalized Experiments with Participant Generation
# Renames the 5 task-order sheets and refreshes the generated stim/file-name
# values (column B) to match the newly-generated participant run.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---------------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Name = "GNG_TO-16502912733290226"
$ws.Range("B2").Value = "go_stims-1650291273274681.csv"
$ws.Range("B3").Value = "GNG_stims-16502912732972128.csv"
$ws.Range("B4").Value = "go_stims-1650291273298216.csv"
$ws.Range("B5").Value = "GNG_stims-16502912733280246.csv"

# --- Sheet 2: NB -----------------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Name = "NB_TO-16502912751280103"
$ws.Range("B2").Value = "OB-16502912741268868.csv"
$ws.Range("B3").Value = "ZB-match_6-16502912737643564.csv"
$ws.Range("B4").Value = "OB-16502912741931415.csv"
$ws.Range("B5").Value = "ZB-match_2-16502912740440402.csv"
$ws.Range("B6").Value = "TB-16502912751061578.csv"
$ws.Range("B7").Value = "ZB-match_0-16502912738972123.csv"
$ws.Range("B8").Value = "OB-16502912741622965.csv"
$ws.Range("B9").Value = "TB-16502912745240514.csv"
$ws.Range("B10").Value = "TB-16502912746882966.csv"

# --- Sheet 3: RS -------------------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Name = "RS_TO-16502912751280103"

# --- Sheet 4: TOL ------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Name = "TOL_TO-1650291275176032"
$ws.Range("B2").Value = "MM_stims-16502912751443212.csv"
$ws.Range("B3").Value = "ZM_stims-16502912751321566.csv"
$ws.Range("B4").Value = "MM_stims-1650291275159978.csv"
$ws.Range("B5").Value = "ZM_stims-16502912751453218.csv"
$ws.Range("B6").Value = "MM_stims-16502912751749415.csv"
$ws.Range("B7").Value = "ZM_stims-16502912751609433.csv"

# --- Sheet 5: vSAT -------------------------------------------------------
$ws = $wb.Worksheets.Item(5)
$ws.Name = "vSAT_TO-16502912752577643"
$ws.Range("B2").Value = "SAT_stims-1650291275180902.csv"
$ws.Range("B3").Value = "vSAT_stims-16502912752418585.csv"
$ws.Range("B4").Value = "vSAT_stims-16502912752217908.csv"
$ws.Range("B5").Value = "SAT_stims-1650291275205986.csv"
